$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($addr, $text)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-CellText "D2" "26.231.99"
Set-CellText "E2" "  -0.42%  "
Set-CellText "D3" "1.661.40"
Set-CellText "E3" "  -0.35%  "
Set-CellText "E4" "  -0.69%  "
Set-CellText "D5" "218.51"
Set-CellText "E5" "  -0.41%  "
Set-CellText "E6" "  -2.01%  "
Set-CellText "D7" "1.004"
Set-CellText "E7" "  -0.64%  "
Set-CellText "D8" "0.2642"
Set-CellText "E8" "  -0.78%  "
Set-CellText "D9" "0.06323"
Set-CellText "E9" "  -1.03%  "
Set-CellText "D10" "20.71"
Set-CellText "E10" "  -0.69%  "
Set-CellText "D11" "0.07792"
Set-CellText "E11" "  -0.73%  "
Set-CellText "B12" "Polkadot"
Set-CellText "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText "D12" "4.506"
Set-CellText "E12" "  -1.10%  "
Set-CellText "B13" "WrappedEther"
Set-CellText "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText "D13" "1.613.16"
Set-CellText "E13" "  -3.69%  "
Set-CellText "D14" "1.889.55"
Set-CellText "E14" "  -0.31%  "
Set-CellText "D15" "0.5648"
Set-CellText "E15" "  +1.86%  "
Set-CellText "D16" "0.0₅8066"
Set-CellText "E16" "  -1.59%  "
Set-CellText "E17" "  -1.05%  "
Set-CellText "D18" "26.229.43"
Set-CellText "E18" "  -0.50%  "
Set-CellText "E19" "  -0.74%  "
Set-CellText "D20" "4.725"
Set-CellText "E20" "  +1.11%  "
Set-CellText "D21" "194.37"
Set-CellText "E21" "  +0.15%  "
Set-CellText "E22" "  -0.38%  "
Set-CellText "D23" "6.018"
Set-CellText "E23" "  -0.44%  "
Set-CellText "D25" "146.36"
Set-CellText "E25" "  +0.17%  "
Set-CellText "D26" "0.1211"
Set-CellText "E26" "  -1.25%  "
Set-CellText "D27" "7.242"
Set-CellText "E27" "  +0.35%  "
Set-CellText "D28" "16.06"
Set-CellText "E28" "  -0.31%  "
Set-CellText "D29" "1.485"
Set-CellText "E29" "  -0.99%  "
Set-CellText "D30" "0.05662"
Set-CellText "E30" "  -3.35%  "
Set-CellText "D31" "1.274"
Set-CellText "E31" "  -0.76%  "
Set-CellText "D32" "3.483"
Set-CellText "E32" "  -2.88%  "
Set-CellText "D33" "3.365"
Set-CellText "E33" "  +2.37%  "
Set-CellText "D34" "1.611"
Set-CellText "E34" "  +0.13%  "
Set-CellText "D35" "2.805"
Set-CellText "E35" "  -0.89%  "
Set-CellText "D36" "0.9437"
Set-CellText "E36" "  -2.78%  "
Set-CellText "E37" "  -0.81%  "
Set-CellText "D38" "0.5771"
Set-CellText "E38" "  -0.93%  "
Set-CellText "D39" "0.01604"
Set-CellText "E39" "  -0.29%  "
Set-CellText "D40" "5.989"
Set-CellText "E40" "  +2.57%  "
Set-CellText "D41" "1.069.12"
Set-CellText "E41" "  -0.23%  "
Set-CellText "D42" "2.582"
Set-CellText "E42" "  -0.02%  "
Set-CellText "D43" "0.8492"
Set-CellText "E43" "  -1.75%  "
Set-CellText "D45" "103.12"
Set-CellText "E45" "  -1.06%  "
Set-CellText "D46" "1.800.20"
Set-CellText "E46" "  -0.28%  "
Set-CellText "D47" "58.30"
Set-CellText "E47" "  +0.41%  "
Set-CellText "E48" "  +2.26%  "
Set-CellText "D49" "1.003"
Set-CellText "E49" "  -1.33%  "
Set-CellText "D50" "0.05333"
Set-CellText "E50" "  +3.27%  "
Set-CellText "B51" "Mantle"
Set-CellText "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText "D51" "0.4350"
Set-CellText "E51" "  -1.00%  "
